# Add an "Electrode Locations" column (C) derived from the file name in column A,
# then sort the data rows (2..last) by electrode location in natural order
# (letter prefix, then numeric suffix), e.g. A1, A3, A5, ..., A11, A15, B13, B14, C1, ...

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# ---- 1. Read existing data (columns A and B) for rows 2..lastRow ----
$items = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $fname = $ws.Cells.Item($r, 1).Value2
    $val = $ws.Cells.Item($r, 2).Value2

    $letters = ""
    $digits = 0
    $loc = ""
    if ($fname -match '^([A-Za-z]+)(\d+)_') {
        $letters = $Matches[1]
        $digits = [int]$Matches[2]
        $loc = "$letters$digits"
    }

    $sortKey = "{0}_{1}" -f $letters, ([string]::Format('{0:D6}', $digits))

    $items += [PSCustomObject]@{
        FileName = $fname
        Value    = $val
        Loc      = $loc
        SortKey  = $sortKey
    }
}

# ---- 2. Sort by the natural electrode-location key (A1, A3, ... A11, A15, B13, ...) ----
$sortedItems = $items | Sort-Object SortKey

# ---- 3. Write the header for the new column C, matching the style of A1/B1 ----
$ws.Range("C1").Value = "Electrode Locations"
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---- 4. Write the sorted rows back (columns A, B, C) ----
$i = 0
foreach ($it in $sortedItems) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $it.FileName
    $ws.Cells.Item($r, 2).Value = $it.Value
    $ws.Cells.Item($r, 3).Value = $it.Loc
    $i = $i + 1
}

# ---- 5. Fix up the sheet dimension to reflect the new column ----
$ws.UsedRange | Out-Null
